$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new deliverable "Plan de Iteración" at B4 (after "Estimacion"),
# shifting the existing B4:B7 entries down to B5:B8.
for ($r = 7; $r -ge 4; $r--) {
    $src = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r + 1, 2).Value2 = $src
}
$ws.Cells.Item(4, 2).Value2 = "Plan de Iteración"

# Update the active selection to match the edited workbook state.
$ws.Range("D6").Select()
